# GestionPermisos_CambioLinea.xlsx edit script
# - Adds a new "total_clp" column (H) with computed CLP totals
# - Updates the "turno" quantity values in column G
# - Applies header + banded-row styling (blue accent1 theme) to the new column
# - Updates the active selection

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Update column G values (turno quantity changed from the old duplicated
#    52.961 / 60 placeholders to the real per-row quantities)
# ---------------------------------------------------------------------------
$ws.Range("G2").Value = 6
$ws.Range("G3").Value = 6
$ws.Range("G4").Value = 6
$ws.Range("G5").Value = 6
$ws.Range("G6").Value = 0
$ws.Range("G7").Value = 6
$ws.Range("G8").Value = 6
$ws.Range("G9").Value = 6

# ---------------------------------------------------------------------------
# 2) Add the new "total_clp" header + values in column H
# ---------------------------------------------------------------------------
$ws.Range("H1").Value = "total_clp"
$ws.Range("H2").Value = 15962400
$ws.Range("H3").Value = 14366160
$ws.Range("H4").Value = 11447849
$ws.Range("H5").Value = 5724294
$ws.Range("H6").Value = 0
$ws.Range("H7").Value = 10567700
$ws.Range("H8").Value = 5283850
$ws.Range("H9").Value = 5283850

# ---------------------------------------------------------------------------
# 3) Styling - thin accent1-colored border (right/top/bottom) around every
#    cell in the new column, matching the banded table look. Applied first
#    (while every H cell still shares the same default formatting) so the
#    border is established uniformly before the fill/font differences below.
# ---------------------------------------------------------------------------
$borderOle = 0xd7b395   # ~ theme Accent1, tint 0.4

foreach ($addr in @("H1","H2","H3","H4","H5","H6","H7","H8","H9")) {
    $cell = $ws.Range($addr)
    foreach ($edge in @([Microsoft.Office.Interop.Excel.XlBordersIndex]::xlEdgeRight,
                         [Microsoft.Office.Interop.Excel.XlBordersIndex]::xlEdgeTop,
                         [Microsoft.Office.Interop.Excel.XlBordersIndex]::xlEdgeBottom)) {
        $b = $cell.Borders.Item($edge)
        $b.Color = $borderOle
        $b.LineStyle = [Microsoft.Office.Interop.Excel.XlLineStyle]::xlContinuous
        $b.Weight = [Microsoft.Office.Interop.Excel.XlBorderWeight]::xlThin
    }
}

# ---------------------------------------------------------------------------
# 4) Styling - header cell H1: bold white text on solid accent1 (theme) fill
# ---------------------------------------------------------------------------
$headerCell = $ws.Range("H1")
$headerCell.Font.Bold = $true
$headerCell.Font.ThemeColor = [Microsoft.Office.Interop.Excel.XlThemeColor]::xlThemeColorLight1
$headerCell.Interior.ThemeColor = [Microsoft.Office.Interop.Excel.XlThemeColor]::xlThemeColorAccent1

# ---------------------------------------------------------------------------
# 5) Styling - banded body rows in column H (alternating light accent1 fill)
#    Even data rows (H2,H4,H6,H8) get a light (20%) accent1 fill.
#    Odd data rows (H3,H5,H7,H9) keep no fill.
# ---------------------------------------------------------------------------
$lightFillOle = 0xf2e6dc   # ~ theme Accent1, tint 0.8 (20% Accent1)

foreach ($addr in @("H2","H4","H6","H8")) {
    $ws.Range($addr).Interior.Color = $lightFillOle
}

# ---------------------------------------------------------------------------
# 6) Update the selected cell / active selection
# ---------------------------------------------------------------------------
$ws.Range("J8").Select()
